$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '28.374.09'
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = '  +1.22%  '

# Row 3
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '1.874.20'
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = '  +0.24%  '

# Row 4
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '1.011'
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = '  +0.95%  '

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '315.54'
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = '  +0.99%  '

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '1.012'
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = '  +1.02%  '

# Row 7
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.5129'

# Row 8
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.3917'
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = '  +1.07%  '

# Row 9
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.08286'
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = '  -0.89%  '

# Row 10
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '1.120'
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = '  +0.72%  '

# Row 11
$ws.Range("B11").Value = 'OKB'
$ws.Range("C11").Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '41.87'
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = '  +0.88%  '

# Row 12
$ws.Range("B12").Value = 'Polkadot'
$ws.Range("C12").Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '6.258'
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = '  +1.44%  '

# Row 13
$ws.Range("B13").Value = 'WrappedEther'
$ws.Range("C13").Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '1.873.39'
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = '  -0.11%  '

# Row 14
$ws.Range("B14").Value = 'Solana'
$ws.Range("C14").Value = 'https://coinranking.com/coin/zNZHO_Sjf+solana-sol'
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '20.34'
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = '  -0.71%  '

# Row 15
$ws.Range("B15").Value = 'Chainlink'
$ws.Range("C15").Value = 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '7.243'
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = '  -0.38%  '

# Row 16
$ws.Range("B16").Value = 'BinanceUSD'
$ws.Range("C16").Value = 'https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd'
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '1.012'
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = '  +1.04%  '

# Row 17
$ws.Range("B17").Value = 'ShibaInu'
$ws.Range("C17").Value = 'https://coinranking.com/coin/xz24e0BjL+shibainu-shib'
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '0.00001103'
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = '  +0.40%  '

# Row 18
$ws.Range("B18").Value = 'Litecoin'
$ws.Range("C18").Value = 'https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc'
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '91.27'
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = '  +0.44%  '

# Row 19
$ws.Range("B19").Value = 'TRON'
$ws.Range("C19").Value = 'https://coinranking.com/coin/qUhEFk1I61atv+tron-trx'
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '0.06716'
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = '  +1.16%  '

# Row 20
$ws.Range("B20").Value = 'Avalanche'
$ws.Range("C20").Value = 'https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax'
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '17.71'
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = '  +0.43%  '

# Row 21
$ws.Range("B21").Value = 'Dai'
$ws.Range("C21").Value = 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '1.013'
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = '  +1.06%  '

# Row 22
$ws.Range("B22").Value = 'Uniswap'
$ws.Range("C22").Value = 'https://coinranking.com/coin/_H5FVG9iW+uniswap-uni'
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '5.986'
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = '  -0.39%  '

# Row 23
$ws.Range("B23").Value = 'WrappedBTC'
$ws.Range("C23").Value = 'https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc'
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '28.421.17'
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = '  +1.27%  '

# Row 24
$ws.Range("B24").Value = 'Cosmos'
$ws.Range("C24").Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '11.15'
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = '  +0.73%  '

# Row 25
$ws.Range("B25").Value = 'Toncoin'
$ws.Range("C25").Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '2.263'
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = '  +0.75%  '

# Row 26
$ws.Range("B26").Value = 'WrappedliquidstakedEther2.0'
$ws.Range("C26").Value = 'https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth'
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '2.106.35'
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = '  +1.18%  '

# Row 27
$ws.Range("B27").Value = 'Monero'
$ws.Range("C27").Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '160.84'
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = '  +1.71%  '

# Row 28
$ws.Range("B28").Value = 'EthereumClassic'
$ws.Range("C28").Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '20.64'
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = '  +0.70%  '

# Row 29
$ws.Range("B29").Value = 'LidoDAOToken'
$ws.Range("C29").Value = 'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo'
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '2.434'
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = '  -1.50%  '

# Row 30
$ws.Range("B30").Value = 'BitcoinCash'
$ws.Range("C30").Value = 'https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch'
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '125.99'
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = '  +1.07%  '

# Row 31
$ws.Range("B31").Value = 'Stellar'
$ws.Range("C31").Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '0.1065'
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = '  +0.37%  '

# Row 32
$ws.Range("B32").Value = 'ImmutableX'
$ws.Range("C32").Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '1.043'
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = '  +1.13%  '

# Row 33
$ws.Range("B33").Value = 'Filecoin'
$ws.Range("C33").Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '5.890'
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = '  +0.30%  '

# Row 34
$ws.Range("B34").Value = 'HuobiToken'
$ws.Range("C34").Value = 'https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht'
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '3.634'
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = '  +1.12%  '

# Row 35
$ws.Range("B35").Value = 'VeChain'
$ws.Range("C35").Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '0.02444'
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = '  +0.53%  '

# Row 36
$ws.Range("B36").Value = 'Hedera'
$ws.Range("C36").Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.06504'
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = '  -0.41%  '

# Row 37
$ws.Range("B37").Value = 'FraxShare'
$ws.Range("C37").Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '9.208'
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = '  -2.18%  '

# Row 38
$ws.Range("B38").Value = 'Algorand'
$ws.Range("C38").Value = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.2184'
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = '  +0.07%  '

# Row 39
$ws.Range("B39").Value = 'TrustWalletToken'
$ws.Range("C39").Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '1.261'
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = '  +3.19%  '

# Row 40
$ws.Range("B40").Value = 'TheSandbox'
$ws.Range("C40").Value = 'https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand'
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.6472'
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = '  +0.04%  '

# Row 41
$ws.Range("B41").Value = 'ARBITRUM'
$ws.Range("C41").Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '1.186'
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = '  -1.15%  '

# Row 42
$ws.Range("B42").Value = 'InternetComputer(DFINITY)'
$ws.Range("C42").Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '4.973'
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = '  -0.39%  '

# Row 43
$ws.Range("B43").Value = 'Aptos'
$ws.Range("C43").Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '11.15'
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = '  -1.44%  '

# Row 44
$ws.Range("B44").Value = 'Decentraland'
$ws.Range("C44").Value = 'https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana'
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.6051'
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = '  -0.37%  '

# Row 45
$ws.Range("B45").Value = 'EnergySwap'
$ws.Range("C45").Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '13.11'
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = '  +1.05%  '

# Row 46
$ws.Range("B46").Value = 'WEMIXTOKEN'
$ws.Range("C46").Value = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '1.291'
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = '  +0.88%  '

# Row 47
$ws.Range("B47").Value = 'PancakeSwap'
$ws.Range("C47").Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '3.694'
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = '  +0.67%  '

# Row 48
$ws.Range("B48").Value = 'NEARProtocol'
$ws.Range("C48").Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '2.022'
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = '  +0.90%  '

# Row 49
$ws.Range("B49").Value = 'EOS'
$ws.Range("C49").Value = 'https://coinranking.com/coin/iAzbfXiBBKkR6+eos-eos'
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '1.217'
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = '  +0.19%  '

# Row 50
$ws.Range("B50").Value = 'Quant'
$ws.Range("C50").Value = 'https://coinranking.com/coin/bauj_21eYVwso+quant-qnt'
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '122.00'
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = '  +0.81%  '

# Row 51
$ws.Range("B51").Value = 'Cronos'
$ws.Range("C51").Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.06898'
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = '  +0.62%  '
